# Auto-generated Excel COM-interop edit script
# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.271.17'
$ws.Range('D3').Value = '1.564.13'
$ws.Range('E3').Value = '  -3.61%  '
$ws.Range('E4').Value = '  -0.42%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '207.32'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -3.08%  '
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('E7').Value = '  -5.58%  '
$ws.Range('E8').Value = '  -2.51%  '
$ws.Range('E9').Value = '  -2.14%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '17.76'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('D12').Value = '1.782.02'
$ws.Range('E12').Value = '  -3.56%  '
$ws.Range('D13').Value = '1.560.41'
$ws.Range('E13').Value = '  -4.94%  '
$ws.Range('E14').Value = '  -3.47%  '
$ws.Range('E15').Value = '  -3.23%  '
$ws.Range('D16').Value = '25.257.60'
$ws.Range('E16').Value = '  -2.68%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '59.31'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -3.01%  '
$ws.Range('E18').Value = '  -3.60%  '
$ws.Range('E19').Value = '  -0.39%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '185.61'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -1.93%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '4.12'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -2.62%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '9.27'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -2.81%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '5.86'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -3.17%  '
$ws.Range('E24').Value = '  -2.15%  '
$ws.Range('E25').Value = '  -0.45%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '139.92'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -2.60%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '1.63'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -7.80%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '14.85'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -1.82%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '6.44'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -3.74%  '
$ws.Range('E30').Value = '  -6.16%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.0463'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -3.73%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '3.04'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -2.88%  '
$ws.Range('E33').Value = '  -3.87%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.46'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -1.51%  '
$ws.Range('E35').Value = '  -4.12%  '
$ws.Range('D36').Value = '1.087.03'
$ws.Range('E36').Value = '  -3.46%  '
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('E38').Value = '  -4.66%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.0149'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -2.09%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.494'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -3.87%  '
$ws.Range('E41').Value = '  -8.87%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.763'
$c.Style = "Normal"
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '93.00'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -4.91%  '
$ws.Range('E44').Value = '  -2.44%  '
$ws.Range('D45').Value = '1.696.47'
$ws.Range('E45').Value = '  -3.54%  '
$ws.Range('D46').Value = '0.0₆0107'
$ws.Range('E46').Value = '  -4.57%  '
$ws.Range('E47').Value = '  -3.07%  '
$ws.Range('E48').Value = '  -3.86%  '
$ws.Range('E49').Value = '  -3.98%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.406'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -1.84%  '
$ws.Range('E51').Value = '  -0.66%  '
